$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (prices). Excel auto-converts
# plain numeric strings assigned via .Value into real numbers (losing
# formatting like trailing zeros / thousands-dot grouping), so those cells
# are pre-formatted as Text ("@") before the value is written, exactly as
# the source data keeps them as literal strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.972.56"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.331.69"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.96"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.84"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  -1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.37"
$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.10"
$ws.Range("E11").Value = "  +2.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0786"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.121"
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  -2.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.686.78"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.330.44"
$ws.Range("E16").Value = "  -1.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.925.50"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.27"
$ws.Range("E19").Value = "  -3.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0892"
$ws.Range("E21").Value = "  -0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.99"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("E23").Value = "  +4.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.45"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.54"
$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.32"
$ws.Range("E30").Value = "  -2.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "145.87"
$ws.Range("E31").Value = "  -11.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.02"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.84"
$ws.Range("E34").Value = "  -2.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0703"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.42"
$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  +2.99%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("E38").Value = "  -1.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("E39").Value = "  -0.85%  "

$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.00"
$ws.Range("E41").Value = "  +21.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("E42").Value = "  -0.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.928.44"
$ws.Range("E43").Value = "  -3.53%  "

$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.17"
$ws.Range("E45").Value = "  -2.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").Value = "  -1.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.555.27"
$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.62"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.87"
$ws.Range("E51").Value = "  +1.49%  "
